# Fixed bug in array expansion
#
# Adds the lookup table (rows 15-19, columns C:F) and an array-entered
# INDEX/MATCH formula (rows 22-25, column D) to the "Referencing" sheet,
# then moves the active selection to C23 to match the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Referencing")

# Header / x-values row
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 4

# Lookup table body (rows 16-19)
$ws.Range("C16").Value = 1.4535833325868115
$ws.Range("D16").Value = 1.4535833325868115
$ws.Range("E16").Value = 1.5117266658902839
$ws.Range("F16").Value = 1.5407983325420203

$ws.Range("C17").Value = 9.0545454545454547
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 18
$ws.Range("F17").Value = 18

$ws.Range("C18").Value = 0.36811506356713858
$ws.Range("D18").Value = 0.36811506356713858
$ws.Range("E18").Value = 0.40588480110308967
$ws.Range("F18").Value = 0.42190146532760275

$ws.Range("C19").Value = 0.65100000000000002
$ws.Range("D19").Value = 0.65100000000000002
$ws.Range("E19").Value = 0.65100000000000002
$ws.Range("F19").Value = 0.65100000000000002

# Lookup key and the (fixed) array formula that now correctly expands
# over D22:D25 instead of collapsing to a single cell.
$ws.Range("C22").Value = 4
$ws.Range("D22:D25").FormulaArray = "=INDEX(C16:F19,,MATCH(C22,C15:F15,0))"

# Match the workbook's saved selection state.
$ws.Range("C23").Select()
